# chore: update Sheets via scheduled runner
# Refreshes market-price derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) in columns H-N across several rows on the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 45000
$ws.Range("J26").Value = 40000
$ws.Range("L26").Value = 40000
$ws.Range("N26").Value = -40688

$ws.Range("H32").Value = 3334.2856
$ws.Range("J32").Value = 3334.2856
$ws.Range("L32").Value = 3334.2856
$ws.Range("N32").Value = -3986.2856

$ws.Range("H33").Value = 171.41176
$ws.Range("I33").Value = 116.35714
$ws.Range("J33").Value = 428.33334
$ws.Range("K33").Value = 116.35714
$ws.Range("L33").Value = 428.33334
$ws.Range("M33").Value = 112.64286
$ws.Range("N33").Value = -886.33334

$ws.Range("H74").Value = 3255.7144
$ws.Range("I74").Value = 3181.6667
$ws.Range("J74").Value = 3700
$ws.Range("K74").Value = 3181.6667
$ws.Range("L74").Value = 3700
$ws.Range("M74").Value = -2245.6667
$ws.Range("N74").Value = -5572

$ws.Range("H77").Value = 3255.7144
$ws.Range("I77").Value = 3181.6667
$ws.Range("J77").Value = 3700
$ws.Range("K77").Value = 15908.3335
$ws.Range("L77").Value = 18500
$ws.Range("M77").Value = -11228.3335
$ws.Range("N77").Value = -27860

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 204473.23
$ws.Range("I132").Value = 209639.9
$ws.Range("J132").Value = 157973.14
$ws.Range("K132").Value = 628919.7
$ws.Range("L132").Value = 473919.42
$ws.Range("M132").Value = -626389.7
$ws.Range("N132").Value = -478979.42

$ws.Range("H137").Value = 15385471
$ws.Range("I137").Value = 18519066
$ws.Range("J137").Value = 2365.2727
$ws.Range("K137").Value = 55557198
$ws.Range("L137").Value = 7095.8181
$ws.Range("M137").Value = -55554648
$ws.Range("N137").Value = -12195.8181

$ws.Range("H138").Value = 1104.2122
$ws.Range("I138").Value = 485
$ws.Range("J138").Value = 1663.8846
$ws.Range("K138").Value = 1455
$ws.Range("L138").Value = 4991.6538
$ws.Range("M138").Value = 3685
$ws.Range("N138").Value = -15271.6538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24882.08
$ws.Range("I32").Value = 4603.6816
$ws.Range("J32").Value = 173590.33
$ws.Range("K32").Value = 4603.6816
$ws.Range("L32").Value = 173590.33
$ws.Range("M32").Value = -4316.6816
$ws.Range("N32").Value = -174164.33

$ws.Range("H63").Value = 32894.75
$ws.Range("I63").Value = 32894.75
$ws.Range("K63").Value = 32894.75
$ws.Range("M63").Value = -32208.75

$ws.Range("H66").Value = 32894.75
$ws.Range("I66").Value = 32894.75
$ws.Range("K66").Value = 164473.75
$ws.Range("M66").Value = -161041.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 6669867.5
$ws.Range("I7").Value = 4801.5
$ws.Range("J7").Value = 20000000
$ws.Range("K7").Value = 4801.5
$ws.Range("L7").Value = 20000000
$ws.Range("M7").Value = -4688.5
$ws.Range("N7").Value = -20000226

$ws.Range("H22").Value = 199.25
$ws.Range("I22").Value = 199
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 199
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -26
$ws.Range("N22").Value = -546

$ws.Range("H134").Value = 14086424
$ws.Range("I134").Value = 18520168
$ws.Range("J134").Value = 2764.5881
$ws.Range("K134").Value = 55560504
$ws.Range("L134").Value = 8293.764299999999
$ws.Range("M134").Value = -55557969
$ws.Range("N134").Value = -13363.7643

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1434.5352
$ws.Range("I31").Value = 855.0877
$ws.Range("J31").Value = 3793.7144
$ws.Range("K31").Value = 855.0877
$ws.Range("L31").Value = 3793.7144
$ws.Range("M31").Value = -560.0877
$ws.Range("N31").Value = -4383.7144

$ws.Range("H34").Value = 1434.5352
$ws.Range("I34").Value = 855.0877
$ws.Range("J34").Value = 3793.7144
$ws.Range("K34").Value = 855.0877
$ws.Range("L34").Value = 3793.7144
$ws.Range("M34").Value = -653.0877
$ws.Range("N34").Value = -4197.7144

$ws.Range("H132").Value = 1637.0308
$ws.Range("I132").Value = 1450.2115
$ws.Range("J132").Value = 2384.3076
$ws.Range("K132").Value = 4350.6345
$ws.Range("L132").Value = 7152.9228
$ws.Range("M132").Value = -1820.6345
$ws.Range("N132").Value = -12212.9228

$ws.Range("H134").Value = 2160.8367
$ws.Range("I134").Value = 1429.1666
$ws.Range("K134").Value = 4287.4998
$ws.Range("M134").Value = -1752.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 525.4286
$ws.Range("I11").Value = 529.6667
$ws.Range("K11").Value = 1589.0001
$ws.Range("M11").Value = -1449.0001

$ws.Range("H113").Value = 11628542
$ws.Range("I113").Value = 651.8
$ws.Range("J113").Value = 21739750
$ws.Range("K113").Value = 1955.4
$ws.Range("L113").Value = 65219250
$ws.Range("M113").Value = 214.6000000000001
$ws.Range("N113").Value = -65223590

$ws.Range("H121").Value = 702.6667
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 783.2
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 2349.6
$ws.Range("M121").Value = 410
$ws.Range("N121").Value = -4969.6

$ws.Range("H130").Value = 999.5238000000001
$ws.Range("J130").Value = 999.5238000000001
$ws.Range("L130").Value = 2998.5714
$ws.Range("N130").Value = -13038.5714

$ws.Range("H131").Value = 6537658
$ws.Range("I131").Value = 569.8889
$ws.Range("J131").Value = 7938462.5
$ws.Range("K131").Value = 1709.6667
$ws.Range("L131").Value = 23815387.5
$ws.Range("M131").Value = 3330.3333
$ws.Range("N131").Value = -23825467.5

$ws.Range("H137").Value = 2973851.2
$ws.Range("I137").Value = 4548307
$ws.Range("J137").Value = 87349.5
$ws.Range("K137").Value = 13644921
$ws.Range("L137").Value = 262048.5
$ws.Range("M137").Value = -13639821
$ws.Range("N137").Value = -272248.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3232

$ws.Range("H17").Value = 18669.666
$ws.Range("J17").Value = 3004.5
$ws.Range("L17").Value = 3004.5
$ws.Range("N17").Value = -3340.5

$ws.Range("H80").Value = 2193
$ws.Range("I80").Value = 2022.591
$ws.Range("J80").Value = 2728.5715
$ws.Range("K80").Value = 2022.591
$ws.Range("L80").Value = 2728.5715
$ws.Range("M80").Value = -1024.591
$ws.Range("N80").Value = -4724.5715

$ws.Range("H83").Value = 2193
$ws.Range("I83").Value = 2022.591
$ws.Range("J83").Value = 2728.5715
$ws.Range("K83").Value = 10112.955
$ws.Range("L83").Value = 13642.8575
$ws.Range("M83").Value = -5120.955
$ws.Range("N83").Value = -23626.8575

$ws.Range("H132").Value = 2375.6204
$ws.Range("I132").Value = 2069.3547
$ws.Range("J132").Value = 3492.5881
$ws.Range("K132").Value = 6208.0641
$ws.Range("L132").Value = 10477.7643
$ws.Range("M132").Value = -3678.0641
$ws.Range("N132").Value = -15537.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 50000
$ws.Range("I13").Value = 50000
$ws.Range("K13").Value = 50000
$ws.Range("M13").Value = -49860

$ws.Range("H55").Value = 518.5238000000001
$ws.Range("I55").Value = 674.5
$ws.Range("J55").Value = 481.82352
$ws.Range("K55").Value = 674.5
$ws.Range("L55").Value = 481.82352
$ws.Range("M55").Value = -501.5
$ws.Range("N55").Value = -827.8235199999999

$ws.Range("H100").Value = 3226.3157
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 3286.6667
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 3286.6667
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4368.6667

$ws.Range("H136").Value = 3575.7693
$ws.Range("I136").Value = 2095.8572
$ws.Range("J136").Value = 16525
$ws.Range("K136").Value = 6287.571599999999
$ws.Range("L136").Value = 49575
$ws.Range("M136").Value = -3737.571599999999
$ws.Range("N136").Value = -54675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3592.7666
$ws.Range("I81").Value = 1906.2307
$ws.Range("J81").Value = 4882.4707
$ws.Range("K81").Value = 3812.4614
$ws.Range("L81").Value = 9764.9414
$ws.Range("M81").Value = -2751.4614
$ws.Range("N81").Value = -11886.9414

$ws.Range("H84").Value = 3592.7666
$ws.Range("I84").Value = 1906.2307
$ws.Range("J84").Value = 4882.4707
$ws.Range("K84").Value = 19062.307
$ws.Range("L84").Value = 48824.70699999999
$ws.Range("M84").Value = -13758.307
$ws.Range("N84").Value = -59432.70699999999

$ws.Range("H132").Value = 3579.75
$ws.Range("I132").Value = 3828.75
$ws.Range("J132").Value = 2832.75
$ws.Range("K132").Value = 11486.25
$ws.Range("L132").Value = 8498.25
$ws.Range("M132").Value = -8956.25
$ws.Range("N132").Value = -13558.25

$ws.Range("H136").Value = 30743.854
$ws.Range("I136").Value = 44283.086
$ws.Range("J136").Value = 2434.5454
$ws.Range("K136").Value = 132849.258
$ws.Range("L136").Value = 7303.6362
$ws.Range("M136").Value = -130299.258
$ws.Range("N136").Value = -12403.6362
